$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 339.33334
$ws.Range("I11").Value = 339.33334
$ws.Range("K11").Value = 339.33334
$ws.Range("M11").Value = -199.33334
$ws.Range("H51").Value = 4298.4
$ws.Range("I51").Value = 3831
$ws.Range("J51").Value = 4999.5
$ws.Range("K51").Value = 3831
$ws.Range("L51").Value = 4999.5
$ws.Range("M51").Value = -3347
$ws.Range("N51").Value = -5967.5
$ws.Range("H92").Value = 360
$ws.Range("I92").Value = 360
$ws.Range("K92").Value = 360
$ws.Range("M92").Value = 888
$ws.Range("H95").Value = 30624
$ws.Range("J95").Value = 30624
$ws.Range("L95").Value = 30624
$ws.Range("N95").Value = -36116
$ws.Range("H113").Value = 15937.5
$ws.Range("I113").Value = 19750
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 19750
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -16496
$ws.Range("N113").Value = -11008
$ws.Range("H138").Value = 2751.7
$ws.Range("J138").Value = 1828.6666
$ws.Range("L138").Value = 5485.9998
$ws.Range("N138").Value = -15765.9998
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 19033.666
$ws.Range("I28").Value = 19033.666
$ws.Range("K28").Value = 19033.666
$ws.Range("M28").Value = -18841.666
$ws.Range("H61").Value = 761
$ws.Range("J61").Value = 990
$ws.Range("L61").Value = 990
$ws.Range("N61").Value = -1414
$ws.Range("H76").Value = 60000
$ws.Range("J76").Value = 60000
$ws.Range("L76").Value = 60000
$ws.Range("N76").Value = -60676
$ws.Range("H79").Value = 60000
$ws.Range("J79").Value = 60000
$ws.Range("L79").Value = 60000
$ws.Range("N79").Value = -62340
$ws.Range("H99").Value = 19033.666
$ws.Range("I99").Value = 19033.666
$ws.Range("K99").Value = 19033.666
$ws.Range("M99").Value = -16038.666
$ws.Range("H102").Value = 4699.5
$ws.Range("I102").Value = 2399
$ws.Range("K102").Value = 2399
$ws.Range("M102").Value = -777
$ws.Range("H122").Value = 2751.182
$ws.Range("I122").Value = 2534.3125
$ws.Range("J122").Value = 3329.5
$ws.Range("K122").Value = 7602.9375
$ws.Range("L122").Value = 9988.5
$ws.Range("M122").Value = -5152.9375
$ws.Range("N122").Value = -14888.5
$ws.Range("H132").Value = 1164.3334
$ws.Range("I132").Value = 1164.3334
$ws.Range("K132").Value = 3493.0002
$ws.Range("M132").Value = -963.0001999999999
$ws.Range("H136").Value = 761
$ws.Range("J136").Value = 990
$ws.Range("L136").Value = 2970
$ws.Range("N136").Value = -8070
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6674.75
$ws.Range("I20").Value = 3899.6667
$ws.Range("K20").Value = 3899.6667
$ws.Range("M20").Value = -3652.6667
$ws.Range("H33").Value = 1700
$ws.Range("I33").Value = 1700
$ws.Range("K33").Value = 1700
$ws.Range("M33").Value = -1364
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 975
$ws.Range("I22").Value = 725
$ws.Range("J22").Value = 1225
$ws.Range("K22").Value = 725
$ws.Range("L22").Value = 1225
$ws.Range("M22").Value = -375
$ws.Range("N22").Value = -1925
$ws.Range("H31").Value = 2395
$ws.Range("I31").Value = 2357.3333
$ws.Range("K31").Value = 2357.3333
$ws.Range("M31").Value = -2062.3333
$ws.Range("H34").Value = 2395
$ws.Range("I34").Value = 2357.3333
$ws.Range("K34").Value = 2357.3333
$ws.Range("M34").Value = -2155.3333
$ws.Range("H58").Value = 1447.5
$ws.Range("I58").Value = 1463.3334
$ws.Range("K58").Value = 1463.3334
$ws.Range("M58").Value = -1260.3334
$ws.Range("H99").Value = 8133.1665
$ws.Range("I99").Value = 8133.1665
$ws.Range("K99").Value = 8133.1665
$ws.Range("M99").Value = -6635.1665
$ws.Range("H105").Value = 1783.3334
$ws.Range("I105").Value = 1841.8334
$ws.Range("J105").Value = 1666.3334
$ws.Range("K105").Value = 1841.8334
$ws.Range("L105").Value = 1666.3334
$ws.Range("M105").Value = -94.83339999999998
$ws.Range("N105").Value = -5160.3334
$ws.Range("H108").Value = 8000
$ws.Range("I108").Value = 8000
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 8000
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -4160
$ws.Range("N108").Value = $null
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = $null
$ws.Range("H126").Value = 8133.1665
$ws.Range("I126").Value = 8133.1665
$ws.Range("K126").Value = 24399.4995
$ws.Range("M126").Value = -21929.4995
$ws.Range("H134").Value = 2316.111
$ws.Range("I134").Value = 2010
$ws.Range("J134").Value = 3387.5
$ws.Range("K134").Value = 6030
$ws.Range("L134").Value = 10162.5
$ws.Range("M134").Value = -3495
$ws.Range("N134").Value = -15232.5
$ws.Range("H136").Value = 1447.5
$ws.Range("I136").Value = 1463.3334
$ws.Range("K136").Value = 4390.0002
$ws.Range("M136").Value = -1840.0002
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 101.3125
$ws.Range("I40").Value = 106
$ws.Range("K40").Value = 424
$ws.Range("M40").Value = -355
$ws.Range("H69").Value = 1896.6875
$ws.Range("I69").Value = 695.6667
$ws.Range("J69").Value = 5499.75
$ws.Range("K69").Value = 2087.0001
$ws.Range("L69").Value = 16499.25
$ws.Range("M69").Value = -1276.0001
$ws.Range("N69").Value = -18121.25
$ws.Range("H72").Value = 1896.6875
$ws.Range("I72").Value = 695.6667
$ws.Range("J72").Value = 5499.75
$ws.Range("K72").Value = 6261.0003
$ws.Range("L72").Value = 49497.75
$ws.Range("M72").Value = -2205.0003
$ws.Range("N72").Value = -57609.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 1500
$ws.Range("I133").Value = 1500
$ws.Range("K133").Value = 4500
$ws.Range("M133").Value = 560
$ws.Range("H70").Value = 12000
$ws.Range("I70").Value = 12000
$ws.Range("K70").Value = 12000
$ws.Range("M70").Value = -11730
$ws.Range("H73").Value = 12000
$ws.Range("I73").Value = 12000
$ws.Range("K73").Value = 12000
$ws.Range("M73").Value = -11064
$ws.Range("H113").Value = 1892
$ws.Range("I113").Value = 1790.6666
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1790.6666
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 379.3334
$ws.Range("N113").Value = -6840
$ws.Range("H126").Value = 3613
$ws.Range("I126").Value = 3415.8
$ws.Range("K126").Value = 10247.4
$ws.Range("M126").Value = -7777.400000000001
$ws.Range("H132").Value = 9500
$ws.Range("I132").Value = 9500
$ws.Range("K132").Value = 28500
$ws.Range("M132").Value = -25970
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = $null
$ws.Range("H46").Value = 1788.6154
$ws.Range("I46").Value = 1375
$ws.Range("J46").Value = 2143.1428
$ws.Range("K46").Value = 1375
$ws.Range("L46").Value = 2143.1428
$ws.Range("M46").Value = -1187
$ws.Range("N46").Value = -2519.1428
$ws.Range("H132").Value = 2681.1904
$ws.Range("I132").Value = 2606.625
$ws.Range("J132").Value = 2919.8
$ws.Range("K132").Value = 7819.875
$ws.Range("L132").Value = 8759.400000000001
$ws.Range("M132").Value = -5289.875
$ws.Range("N132").Value = -13819.4
$ws.Range("H136").Value = 3909.4546
$ws.Range("I136").Value = 3444.889
$ws.Range("K136").Value = 10334.667
$ws.Range("M136").Value = -7784.667000000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 14814.728
$ws.Range("J14").Value = 13796.2
$ws.Range("L14").Value = 13796.2
$ws.Range("N14").Value = -14132.2
$ws.Range("H122").Value = 1190.95
$ws.Range("I122").Value = 1211
$ws.Range("K122").Value = 3633
$ws.Range("M122").Value = -1183
$ws.Range("H126").Value = 3899.1428
$ws.Range("I126").Value = 3715.6667
$ws.Range("K126").Value = 11147.0001
$ws.Range("M126").Value = -8677.000100000001
$ws.Range("H132").Value = 1777.4667
$ws.Range("I132").Value = 1566.7307
$ws.Range("K132").Value = 4700.1921
$ws.Range("M132").Value = -2170.1921
